$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet has a yearly data table (2007..2020 in columns D..Q). A new
# "2021" column is being appended in column R, mirroring the formatting of
# the preceding data columns (O is a representative interior column, since
# Q currently carries the "last column" styling that should move to R).
[void]$ws.Range("O4:O14").Copy()
[void]$ws.Range("R4:R14").PasteSpecial(-4122)

# Header year
$ws.Range("R4").Value = 2021

# Data values for 2021, one per oblast / republic row
$ws.Range("R5").Value = 1
$ws.Range("R6").Value = 2.2
$ws.Range("R7").Value = 1.7
$ws.Range("R8").Value = "-"
$ws.Range("R9").Value = 0.3
$ws.Range("R10").Value = 1.1
$ws.Range("R11").Value = "-"
$ws.Range("R12").Value = 0.9
$ws.Range("R13").Value = 0.4
$ws.Range("R14").Value = 0.6

# Move the lingering selection one column further right, matching the
# author's cursor position after extending the table.
[void]$ws.Range("S17").Select()
